# Add a new "z_value" column (E) to the worksheet.
#
# The header cell (E1) picks up the same bold/border/centered style already
# used by the other header cells (B1:D1) by copying the format from D1.
# The data cells (E2:E11) are plain numeric values, matching the rest of
# the unstyled data columns.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clone D1's formatting onto E1, then set its text -- this reuses the
# existing bold-header cell style instead of creating a new one.
$ws.Range("D1").Copy()
$ws.Range("E1").PasteSpecial(-4122)
$ws.Range("E1").Value = "z_value"

# Fill in the computed z_value figures for each row.
$values = @(
    90,
    90,
    90,
    90,
    87.05882352941177,
    83.75,
    71.25,
    71.25,
    70.45454545454545,
    68.33333333333333
)

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 5).Value = $values[$i]
}
